# Rename Sheet1 -> Artificial
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Artificial"

# Update the two chart series so their formulas point at the renamed sheet
# (chart series formulas are not automatically re-pointed by a sheet rename).
$lineChart = $ws.ChartObjects(1).Chart
$lineSeries = $lineChart.SeriesCollection(1)
$lineSeries.Formula = "=SERIES(Artificial!B1,Artificial!A:A,Artificial!B:B,1)"

$barChart = $ws.ChartObjects(2).Chart
$barSeries = $barChart.SeriesCollection(1)
$barSeries.Formula = "=SERIES(Artificial!B1,Artificial!`$E`$25:`$I`$25,Artificial!`$E`$26:`$I`$26,1)"

# Refresh the "Stopa zwrotu" (rate of return) figures in E26:I26
$ws.Range("E26").Value = 0.05033011195100934
$ws.Range("F26").Value = 0.005496015388843034
$ws.Range("G26").Value = 0.03840696244442343
$ws.Range("H26").Value = 0.05022962112514362
$ws.Range("I26").Value = 0.07871462264150941

# Append the new daily NAV history rows 87-93 (copy A86's date style down first)
$ws.Range("A86").Copy() | Out-Null
$ws.Range("A87:A93").PasteSpecial(-4122) | Out-Null

$dates = @(43099, 43100, 43101, 43102, 43103, 43104, 43105)
$vals  = @(105.93, 105.93, 105.93, 105.93, 107.65, 109.2, 109.77)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 87 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $vals[$i]
}
